$wb = $excel.ActiveWorkbook

# --- "buses" sheet: debug PGD chart dataframe ---
$wsBuses = $wb.Worksheets.Item("buses")
$wsBuses.Range("C12").Value = 0

# --- "lines" sheet ---
$wsLines = $wb.Worksheets.Item("lines")
$wsLines.Range("D11").Value = 0

# Highlight the recomputed shared-formula cell on "lines" (bold font, new style)
$wsLines.Range("D12").Font.Bold = $true

# Make "lines" the active sheet/tab (matches tabSelected moving to sheet2)
$wsLines.Activate()
$wsLines.Range("D12").Select() | Out-Null
